$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.3884013333333334
$ws.Range("H2").Value = 1.165204
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.757644
$ws.Range("N2").Value = 2.272932
$ws.Range("O2").Value = 0.02401898721285653
$ws.Range("P2").Value = 0.02518910262217759
$ws.Range("Q2").Value = 0.294269939792
$ws.Range("R2").Value = 2.648429458128
$ws.Range("S2").Value = 0.02401898721285653
$ws.Range("T2").Value = 0.02518910262217759

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.3884013333333334
$ws.Range("H3").Value = 1.165204
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 25.23919433333333
$ws.Range("N3").Value = 75.717583
$ws.Range("O3").Value = 0.8001381730141521
$ws.Range("P3").Value = 0.8391179183936208
$ws.Range("Q3").Value = 9.802936731325779
$ws.Range("R3").Value = 88.22643058193202
$ws.Range("S3").Value = 0.8001381730141521
$ws.Range("T3").Value = 0.8391179183936208

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.3884013333333334
$ws.Range("H4").Value = 1.165204
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.6151326666666667
$ws.Range("N4").Value = 1.845398
$ws.Range("O4").Value = 0.0195010633686494
$ws.Range("P4").Value = 0.02045108239083319
$ws.Range("Q4").Value = 0.2389183479102222
$ws.Range("R4").Value = 2.150265131192
$ws.Range("S4").Value = 0.0195010633686494
$ws.Range("T4").Value = 0.02045108239083319

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.3884013333333334
$ws.Range("H5").Value = 1.165204
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.5356743333333333
$ws.Range("N5").Value = 1.607023
$ws.Range("O5").Value = 0.01698205880675987
$ws.Range("P5").Value = 0.01780936132853939
$ws.Range("Q5").Value = 0.2080566252991111
$ws.Range("R5").Value = 1.872509627692
$ws.Range("S5").Value = 0.01698205880675987
$ws.Range("T5").Value = 0.01780936132853939

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.3884013333333334
$ws.Range("H6").Value = 1.165204
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 4.395899500000001
$ws.Range("N6").Value = 8.791799000000001
$ws.Range("O6").Value = 0.139359717597582
$ws.Range("P6").Value = 0.09743253526482902
$ws.Range("Q6").Value = 1.707373226999334
$ws.Range("R6").Value = 10.244239361996
$ws.Range("S6").Value = 0.139359717597582
$ws.Range("T6").Value = 0.09743253526482902
